$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 195.76923
$ws.Range("I18").Value = 195.41667
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 195.41667
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = 88.58332999999999
$ws.Range("N18").Value = -768

$ws.Range("H94").Value = 5246.923
$ws.Range("I94").Value = 3356.6667
$ws.Range("J94").Value = 9500
$ws.Range("K94").Value = 3356.6667
$ws.Range("L94").Value = 9500
$ws.Range("M94").Value = -2905.6667
$ws.Range("N94").Value = -10402

$ws.Range("H121").Value = 1108.3334
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1108.3334
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3325.0002
$ws.Range("N121").Value = -6819.0002

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 10159.923
$ws.Range("I131").Value = 415.8
$ws.Range("J131").Value = 16250
$ws.Range("K131").Value = 1247.4
$ws.Range("L131").Value = 48750
$ws.Range("M131").Value = 3792.6
$ws.Range("N131").Value = -58830

$ws.Range("H132").Value = 171872.62
$ws.Range("I132").Value = 2158.261
$ws.Range("J132").Value = 772400.4
$ws.Range("K132").Value = 6474.782999999999
$ws.Range("L132").Value = 2317201.2
$ws.Range("M132").Value = -3944.782999999999
$ws.Range("N132").Value = -2322261.2

$ws.Range("H138").Value = 1533.85
$ws.Range("I138").Value = 704.6486
$ws.Range("J138").Value = 2020.8413
$ws.Range("K138").Value = 2113.9458
$ws.Range("L138").Value = 6062.5239
$ws.Range("M138").Value = 3026.0542
$ws.Range("N138").Value = -16342.5239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 145715.19
$ws.Range("I32").Value = 166780.98
$ws.Range("J32").Value = 65401.875
$ws.Range("K32").Value = 166780.98
$ws.Range("L32").Value = 65401.875
$ws.Range("M32").Value = -166493.98
$ws.Range("N32").Value = -65975.875

$ws.Range("H45").Value = 1005.55554
$ws.Range("I45").Value = 780
$ws.Range("J45").Value = 1287.5
$ws.Range("K45").Value = 780
$ws.Range("L45").Value = 1287.5
$ws.Range("M45").Value = -403
$ws.Range("N45").Value = -2041.5

$ws.Range("H119").Value = 50000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 50000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676

$ws.Range("H122").Value = 1209.7273
$ws.Range("I122").Value = 1209.7273
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3629.1819
$ws.Range("L122").Value = $null
$ws.Range("M122").Value = -1179.1819
$ws.Range("N122").Value = $null

$ws.Range("H123").Value = 48000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 48000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800

$ws.Range("H132").Value = 2609036.8
$ws.Range("I132").Value = 3189346
$ws.Range("J132").Value = 920864.75
$ws.Range("K132").Value = 9568038
$ws.Range("L132").Value = 2762594.25
$ws.Range("M132").Value = -9565508
$ws.Range("N132").Value = -2767654.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = $null

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 56640
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 56640
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 56640
$ws.Range("N100").Value = -58804

$ws.Range("H103").Value = 13496.714
$ws.Range("I103").Value = 9180
$ws.Range("J103").Value = 24288.5
$ws.Range("K103").Value = 9180
$ws.Range("L103").Value = 24288.5
$ws.Range("M103").Value = -8008
$ws.Range("N103").Value = -26632.5

$ws.Range("H134").Value = 1486.9375
$ws.Range("I134").Value = 1486.9375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4460.8125
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1925.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13226.125
$ws.Range("I5").Value = 849.75
$ws.Range("J5").Value = 25602.5
$ws.Range("K5").Value = 2549.25
$ws.Range("L5").Value = 76807.5
$ws.Range("M5").Value = -2437.25
$ws.Range("N5").Value = -77031.5

$ws.Range("H135").Value = 13226.125
$ws.Range("I135").Value = 849.75
$ws.Range("J135").Value = 25602.5
$ws.Range("K135").Value = 7647.75
$ws.Range("L135").Value = 230422.5
$ws.Range("M135").Value = -5112.75
$ws.Range("N135").Value = -235492.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 30113
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 30113
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 30113
$ws.Range("N51").Value = -31131

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws.Range("H113").Value = 3433.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3433.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3433.3333
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -7773.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5112.5
$ws.Range("I7").Value = 5112.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5112.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5000.5
$ws.Range("N7").Value = $null

$ws.Range("H122").Value = 2464.65
$ws.Range("I122").Value = 2009.0333
$ws.Range("J122").Value = 3831.5
$ws.Range("K122").Value = 6027.0999
$ws.Range("L122").Value = 11494.5
$ws.Range("M122").Value = -3577.0999
$ws.Range("N122").Value = -16394.5

$ws.Range("H126").Value = 5112.5
$ws.Range("I126").Value = 5112.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15337.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12867.5
$ws.Range("N126").Value = $null

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws.Range("H132").Value = 465523.5
$ws.Range("I132").Value = 103710.35
$ws.Range("J132").Value = 1671567.4
$ws.Range("K132").Value = 311131.05
$ws.Range("L132").Value = 5014702.199999999
$ws.Range("M132").Value = -308601.05
$ws.Range("N132").Value = -5019762.199999999

$ws.Range("H136").Value = 627850.75
$ws.Range("I136").Value = 1430983.9
$ws.Range("J136").Value = 3191.6667
$ws.Range("K136").Value = 4292951.699999999
$ws.Range("L136").Value = 9575.000100000001
$ws.Range("M136").Value = -4290401.699999999
$ws.Range("N136").Value = -14675.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 24902.334
$ws.Range("I37").Value = 1026
$ws.Range("J37").Value = 29677.6
$ws.Range("K37").Value = 1026
$ws.Range("L37").Value = 29677.6
$ws.Range("M37").Value = -823
$ws.Range("N37").Value = -30083.6

$ws.Range("H122").Value = 10061.23
$ws.Range("I122").Value = 9004
$ws.Range("K122").Value = 27012
$ws.Range("M122").Value = -24562

$ws.Range("H132").Value = 4944.1035
$ws.Range("I132").Value = 993.7368
$ws.Range("J132").Value = 12449.8
$ws.Range("K132").Value = 2981.2104
$ws.Range("L132").Value = 37349.39999999999
$ws.Range("M132").Value = -451.2103999999999
$ws.Range("N132").Value = -42409.39999999999

$ws.Range("H136").Value = 372014.4
$ws.Range("I136").Value = 1571.9375
$ws.Range("K136").Value = 4715.8125
$ws.Range("M136").Value = -2165.8125
